# overall.xlsx ("Tabelle1"): bump the 10 GB / Xeon E-2276G benchmark figure
# from 53.5 to 53.6 (small data correction) and leave the sheet scrolled /
# selected the way the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Xeon E-2276G row (row 4), "10 GB" column (G): 53.5 -> 53.6
$ws.Range("G4").Value = 53.6

# Restore the view state (scroll position + active selection) the workbook
# was left in after the edit.
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G5").Select()
